# The sheet had two "summary" rows (5 and 6) that each mixed a month
# label together with placeholder values ('.', '.', '0,00'). The fix
# splits each of those into two rows: a standalone month-label row,
# followed by the original month + placeholder-values row, now shifted
# down two rows (5->7, 6->8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing row 5 ("06/2012", ".", ".", "0,00") down to row 7,
# and the existing row 6 ("07/2012", ".", ".", "0,00") down to row 8.
$ws.Range("A7").Value = $ws.Range("A5").Value2
$ws.Range("B7").Value = $ws.Range("B5").Value2
$ws.Range("C7").Value = $ws.Range("C5").Value2
$ws.Range("D7").Value = $ws.Range("D5").Value2

$ws.Range("A8").Value = $ws.Range("A6").Value2
$ws.Range("B8").Value = $ws.Range("B6").Value2
$ws.Range("C8").Value = $ws.Range("C6").Value2
$ws.Range("D8").Value = $ws.Range("D6").Value2

# Clear the now-duplicated B:D values out of rows 5 and 6, leaving
# just a single month-label cell in column A for each.
$ws.Range("B5:D6").ClearContents()

$ws.Range("A5").Value = "06/2012"
$ws.Range("A6").Value = "07/2012"
